# Set level timer to 90s
# - Food_Define: add a "Waiting" column (header H5 + values H6:H17=10), becomes active sheet, selection K13
# - Level_define: time_sec (D7:D10) 120 -> 90, selection E19
# - Incident_Define: loses tabSelected, selection J9

$wb = $excel.ActiveWorkbook

$foodDefine  = $wb.Worksheets.Item("Food_Define")
$levelDefine = $wb.Worksheets.Item("Level_define")
$incidentDef = $wb.Worksheets.Item("Incident_Define")

# --- Food_Define: new "Waiting" column (H) ---------------------------------

# H5 header cell: copy the header formatting used by the rest of row 5 (G5)
# then give it the "Waiting" label.
$foodDefine.Range("G5").Copy()
$foodDefine.Range("H5").PasteSpecial(-4122)  # xlPasteFormats
$foodDefine.Range("H5").Value = "Waiting"

# H6 (the blank banner row under the header): match the rest of row 6 (G6).
$foodDefine.Range("G6").Copy()
$foodDefine.Range("H6").PasteSpecial(-4122)  # xlPasteFormats

# H7:H17 data values - everybody waits 10 "ticks".
$foodDefine.Range("H7:H17").Value = 10

$foodDefine.Application.CutCopyMode = $false

# --- Level_define: level timer 120s -> 90s ----------------------------------

$levelDefine.Range("D7:D10").Value = 90

# --- Selection / active sheet bookkeeping -----------------------------------
# Originally "Incident_Define" was the active tab; now "Food_Define" is.

$incidentDef.Activate()
$incidentDef.Range("J9").Select()

$levelDefine.Activate()
$levelDefine.Range("E19").Select()

$foodDefine.Activate()
$foodDefine.Range("K13").Select()
